$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text, matching the source data which uses
# locale-formatted strings (e.g. thousands-dot grouping) that are not valid numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '24.431.21'
$ws.Range('E2').Value = '  +1.33%  '
$ws.Range('D3').Value = '1.666.22'
$ws.Range('E3').Value = '  +1.41%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').Value = '313.56'
$ws.Range('E5').Value = '  +1.94%  '
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  -0.17%  '
$ws.Range('D7').Value = '0.3967'
$ws.Range('E7').Value = '  +1.45%  '
$ws.Range('D8').Value = '0.3917'
$ws.Range('E8').Value = '  +1.52%  '
$ws.Range('D9').Value = '52.13'
$ws.Range('E9').Value = '  +7.12%  '
$ws.Range('D10').Value = '1.401'
$ws.Range('E10').Value = '  +3.19%  '
$ws.Range('D11').Value = '1.001'
$ws.Range('E11').Value = '  -0.25%  '
$ws.Range('D12').Value = '0.08607'
$ws.Range('E12').Value = '  +1.79%  '
$ws.Range('D13').Value = '24.42'
$ws.Range('E13').Value = '  +1.56%  '
$ws.Range('D14').Value = '7.330'
$ws.Range('E14').Value = '  +2.57%  '
$ws.Range('D15').Value = '0.00001363'
$ws.Range('E15').Value = '  +6.13%  '
$ws.Range('D16').Value = '7.892'
$ws.Range('E16').Value = '  +5.38%  '
$ws.Range('D17').Value = '1.666.20'
$ws.Range('E17').Value = '  +1.22%  '
$ws.Range('D18').Value = '95.42'
$ws.Range('E18').Value = '  +1.12%  '
$ws.Range('D19').Value = '0.06990'
$ws.Range('E19').Value = '  +0.62%  '
$ws.Range('D20').Value = '20.59'
$ws.Range('E20').Value = '  -1.75%  '
$ws.Range('D21').Value = '7.011'
$ws.Range('E21').Value = '  +0.80%  '
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  -0.26%  '
$ws.Range('D23').Value = '13.76'
$ws.Range('E23').Value = '  +0.19%  '
$ws.Range('D24').Value = '24.429.73'
$ws.Range('E24').Value = '  +1.32%  '
$ws.Range('D25').Value = '2.427'
$ws.Range('E25').Value = '  +3.42%  '
$ws.Range('D26').Value = '3.039'
$ws.Range('E26').Value = '  +11.61%  '
$ws.Range('D27').Value = '22.54'
$ws.Range('E27').Value = '  +0.15%  '
$ws.Range('D28').Value = '157.65'
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('D29').Value = '142.87'
$ws.Range('E29').Value = '  +0.85%  '
$ws.Range('D30').Value = '5.457'
$ws.Range('E30').Value = '  +1.65%  '
$ws.Range('D31').Value = '8.104'
$ws.Range('E31').Value = '  -9.41%  '
$ws.Range('D32').Value = '2.514'
$ws.Range('E32').Value = '  +1.60%  '
$ws.Range('D33').Value = '1.847.57'
$ws.Range('E33').Value = '  +1.13%  '
$ws.Range('D34').Value = '1.070'
$ws.Range('E34').Value = '  +9.11%  '
$ws.Range('D35').Value = '0.08286'
$ws.Range('E35').Value = '  +3.17%  '
$ws.Range('D36').Value = '0.03033'
$ws.Range('E36').Value = '  +3.32%  '
$ws.Range('D37').Value = '6.915'
$ws.Range('E37').Value = '  -4.01%  '
$ws.Range('D38').Value = '0.2776'
$ws.Range('E38').Value = '  +2.45%  '
$ws.Range('D39').Value = '11.13'
$ws.Range('E39').Value = '  +10.52%  '
$ws.Range('D40').Value = '0.09252'
$ws.Range('E40').Value = '  +0.09%  '
$ws.Range('D41').Value = '0.7749'
$ws.Range('E41').Value = '  +1.46%  '
$ws.Range('D42').Value = '13.88'
$ws.Range('E42').Value = '  +5.90%  '
$ws.Range('D43').Value = '1.443'
$ws.Range('E43').Value = '  -2.23%  '
$ws.Range('D44').Value = '16.57'
$ws.Range('E44').Value = '  +3.47%  '
$ws.Range('D45').Value = '0.7134'
$ws.Range('E45').Value = '  +3.60%  '
$ws.Range('D46').Value = '2.540'
$ws.Range('E46').Value = '  +1.99%  '
$ws.Range('D47').Value = '4.141'
$ws.Range('E47').Value = '  +1.15%  '
$ws.Range('D48').Value = '1.001'
$ws.Range('E48').Value = '  -0.25%  '
$ws.Range('D49').Value = '0.08463'
$ws.Range('E49').Value = '  +0.63%  '
$ws.Range('D50').Value = '136.70'
$ws.Range('E50').Value = '  +1.90%  '
$ws.Range('D51').Value = '1.275'
$ws.Range('E51').Value = '  +1.01%  '
